$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) cells being updated keep text formatting, otherwise Excel
# would reinterpret strings such as "1.00" or "0.160" as numbers and drop the
# trailing zeros that the site-scraped values always carry.
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D10").NumberFormat = "@"
$ws.Range("D12:D26").NumberFormat = "@"
$ws.Range("D28:D29").NumberFormat = "@"
$ws.Range("D31:D51").NumberFormat = "@"

$ws.Range("D2").Value = '59.289.21'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '2.520.17'
$ws.Range("E3").Value = '  -0.51%  '
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").Value = '535.08'
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("D6").Value = '138.97'
$ws.Range("E6").Value = '  -2.81%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.38%  '
$ws.Range("D8").Value = '0.562'
$ws.Range("E8").Value = '  -1.60%  '
$ws.Range("D9").Value = '2.524.77'
$ws.Range("D10").Value = '0.0991'
$ws.Range("E10").Value = '  -0.17%  '
$ws.Range("E11").Value = '  +1.43%  '
$ws.Range("D12").Value = '5.39'
$ws.Range("E12").Value = '  -1.48%  '
$ws.Range("D13").Value = '0.355'
$ws.Range("E13").Value = '  +1.26%  '
$ws.Range("D14").Value = '2.969.85'
$ws.Range("E14").Value = '  +0.96%  '
$ws.Range("D15").Value = '23.03'
$ws.Range("E15").Value = '  -2.33%  '
$ws.Range("D16").Value = '59.228.37'
$ws.Range("D17").Value = '0.0000140'
$ws.Range("E17").Value = '  +0.70%  '
$ws.Range("D18").Value = '2.516.12'
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("D19").Value = '10.89'
$ws.Range("E19").Value = '  -2.97%  '
$ws.Range("D20").Value = '4.20'
$ws.Range("E20").Value = '  -1.50%  '
$ws.Range("D21").Value = '320.72'
$ws.Range("E21").Value = '  -0.62%  '
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '5.79'
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("D24").Value = '62.15'
$ws.Range("E24").Value = '  +0.28%  '
$ws.Range("D25").Value = '0.421'
$ws.Range("E25").Value = '  -3.62%  '
$ws.Range("D26").Value = '0.165'
$ws.Range("E26").Value = '  +1.33%  '
$ws.Range("E27").Value = '  +0.51%  '
$ws.Range("D28").Value = '7.75'
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("D29").Value = '6.71'
$ws.Range("E29").Value = '  -0.40%  '
$ws.Range("E30").Value = '  +0.47%  '
$ws.Range("D31").Value = '0.0₃0762'
$ws.Range("E31").Value = '  -1.15%  '
$ws.Range("B32").Value = 'USDe'
$ws.Range("C32").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  +0.30%  '
$ws.Range("B33").Value = 'Monero'
$ws.Range("C33").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D33").Value = '159.82'
$ws.Range("E33").Value = '  +1.34%  '
$ws.Range("D34").Value = '1.45'
$ws.Range("E34").Value = '  +1.32%  '
$ws.Range("D35").Value = '1.11'
$ws.Range("E35").Value = '  -6.25%  '
$ws.Range("D36").Value = '18.51'
$ws.Range("E36").Value = '  -0.32%  '
$ws.Range("D37").Value = '4.18'
$ws.Range("E37").Value = '  -3.92%  '
$ws.Range("D38").Value = '1.58'
$ws.Range("E38").Value = '  -1.99%  '
$ws.Range("D39").Value = '36.97'
$ws.Range("E39").Value = '  +0.33%  '
$ws.Range("D40").Value = '3.64'
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("D41").Value = '0.803'
$ws.Range("E41").Value = '  -1.07%  '
$ws.Range("D42").Value = '283.06'
$ws.Range("E42").Value = '  -5.60%  '
$ws.Range("D43").Value = '5.22'
$ws.Range("E43").Value = '  -7.76%  '
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  +0.18%  '
$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D45").Value = '10.86'
$ws.Range("E45").Value = '  +0.85%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = '0.597'
$ws.Range("E46").Value = '  -0.88%  '
$ws.Range("D47").Value = '123.27'
$ws.Range("E47").Value = '  -1.84%  '
$ws.Range("D48").Value = '0.0922'
$ws.Range("D49").Value = '18.53'
$ws.Range("E49").Value = '  -0.61%  '
$ws.Range("D50").Value = '0.0508'
$ws.Range("E50").Value = '  -0.93%  '
$ws.Range("D51").Value = '0.0222'
$ws.Range("E51").Value = '  -2.24%  '
